# Update cryptos list sheet with latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $b, $c, $d, $e) {
    if ($b -ne $null) { $ws.Cells.Item($row, 2).Value = $b }
    if ($c -ne $null) { $ws.Cells.Item($row, 3).Value = $c }
    if ($d -ne $null) { $ws.Cells.Item($row, 4).Value = $d }
    if ($e -ne $null) { $ws.Cells.Item($row, 5).Value = $e }
}

Set-Row 2  $null $null "67.802.96"  "  +0.25%  "
Set-Row 3  $null $null "3.804.84"   "  +0.76%  "
Set-Row 4  $null $null "0.998"      "  -0.19%  "
Set-Row 5  $null $null "604.01"     "  +1.57%  "
Set-Row 6  $null $null "166.06"     "  -0.59%  "
Set-Row 7  $null $null $null        "  +0.05%  "
Set-Row 8  $null $null $null        "  -0.07%  "
Set-Row 9  $null $null $null        "  +0.54%  "
Set-Row 10 $null $null $null        "  +1.00%  "
Set-Row 11 $null $null $null        "  -0.17%  "
Set-Row 12 $null $null "0.0000251"  "  -0.91%  "
Set-Row 13 $null $null "35.92"      "  -0.20%  "
Set-Row 14 $null $null "4.447.73"   "  +0.79%  "
Set-Row 15 $null $null "3.825.66"   "  +1.77%  "
Set-Row 16 $null $null "18.52"      "  +0.92%  "
Set-Row 17 $null $null "67.814.83"  "  +0.32%  "
Set-Row 18 $null $null "7.07"       "  +1.36%  "
Set-Row 19 $null $null $null        "  +1.44%  "
Set-Row 20 $null $null "463.44"     "  +1.65%  "
Set-Row 21 $null $null "9.87"       "  -1.28%  "
Set-Row 22 $null $null "0.702"      "  +1.22%  "
Set-Row 23 $null $null $null        "  -3.56%  "
Set-Row 24 $null $null "83.25"      "  +0.14%  "
Set-Row 25 $null $null "12.14"      "  +2.12%  "
Set-Row 26 $null $null $null        "  -0.30%  "
Set-Row 27 $null $null "10.07"      "  +0.40%  "
Set-Row 28 $null $null $null        "  -0.12%  "
Set-Row 29 $null $null "3.954.91"   "  +0.77%  "
Set-Row 30 $null $null $null        "  +0.99%  "
Set-Row 31 $null $null "7.36"       "  +1.92%  "
Set-Row 32 $null $null $null        "  -0.76%  "
Set-Row 33 $null $null "29.44"      "  -0.54%  "
Set-Row 34 $null $null $null        "  +0.15%  "
Set-Row 35 $null $null $null        "  -0.06%  "
Set-Row 36 $null $null "0.0999"     "  +0.12%  "
Set-Row 37 $null $null $null        "  +0.60%  "

Set-Row 38 "Filecoin" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil" "5.81"  "  +1.37%  "
Set-Row 39 "Mantle"   "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"   "0.995" "  -0.16%  "

Set-Row 40 $null $null $null        "  -3.17%  "
Set-Row 41 $null $null "1.00"       "  +0.03%  "
Set-Row 42 $null $null $null        "  +0.01%  "
Set-Row 43 $null $null "44.63"      "  -1.88%  "
Set-Row 44 $null $null "47.84"      "  -0.85%  "
Set-Row 45 $null $null $null        "  +0.72%  "

Set-Row 46 "Monero"     "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"  "150.59" "  +0.69%  "
Set-Row 47 "ONDO"       "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"       "1.38"   "  +10.89%  "
Set-Row 48 "EnergySwap" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"  "27.53"  "  +7.69%  "

Set-Row 49 $null $null $null        "  +0.64%  "
Set-Row 50 $null $null $null        "  +2.06%  "
Set-Row 51 $null $null "389.96"     "  +0.52%  "
